# Automatische test-sync: 2025-07-23 22:36:50
# Append new mail-log row (#10) on the "Logs" sheet and re-sync the
# "Dashboard" category summary table to reflect the updated counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append the new log entry as row 20 on the "Logs" sheet
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(20, 1).Value = "Kun jij de planning voor volgende week alvast voorbereiden?"
$logs.Cells.Item(20, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(20, 3).Value = "Testmail #10: Kun jij de planning voor volgende week alvast voorbereiden?"
$logs.Cells.Item(20, 4).Value = "Planning / Afspraak"
$logs.Cells.Item(20, 5).Value = "Beste afzender,`nDank je wel voor je e-mail. Ik kan zeker de planning voor volgende week voorbereiden. Zou je wat meer specifieke informatie kunnen geven over de taken en activiteiten die moeten worden ingepland? Zodra ik meer details heb, kan ik aan de slag gaan met het opstellen van de planning.`nMet vriendelijke groet,`n[Naam] [Functie]"
$logs.Cells.Item(20, 6).Value = "2025-07-23 22:36:15"
$logs.Cells.Item(20, 7).Value = "Ja"
$logs.Cells.Item(20, 8).Value = "Nee"
$logs.Cells.Item(20, 9).Value = "Ja"
$logs.Cells.Item(20, 10).Value = "Nee"

# Multi-line content (column E) auto-expands the row height when the
# value is set via COM; restore the default (non-custom) row height so
# the new row matches the rest of the sheet.
$logs.Rows.Item(20).AutoFit()

# ---------------------------------------------------------------------
# 2. Re-sync the "Dashboard" category summary table
#    (Planning / Afspraak now has 2 occurrences and moves up to be
#    right after "Bestelling / Levering"; the remaining single-count
#    categories shift down one row, keeping their relative order.)
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(5, 1).Value = "Planning / Afspraak"
$dash.Cells.Item(5, 2).Value = 2

$dash.Cells.Item(6, 1).Value = "Factuur / Administratie"
$dash.Cells.Item(6, 2).Value = 1

$dash.Cells.Item(7, 1).Value = "IT / Technisch probleem"
$dash.Cells.Item(7, 2).Value = 1

$dash.Cells.Item(8, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(8, 2).Value = 1

# ---------------------------------------------------------------------
# 3. Extend the conditional-formatting ranges on "Logs" to cover the
#    newly appended row 20 (was 2:19, now 2:20 for columns D,G,H,I,J)
# ---------------------------------------------------------------------
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $logs.Range("$col`2:$col`19")
    $newRange = $logs.Range("$col`2:$col`20")
    $conditions = $oldRange.FormatConditions
    for ($i = 1; $i -le $conditions.Count; $i++) {
        $conditions.Item($i).ModifyAppliesToRange($newRange)
    }
}
